$wb = $excel.ActiveWorkbook

# Rename the "Include ValueSets" and "Include from Krebsstadium Cod" sheets
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Include from Krebsstadium Cod").Name = "Include #1"

$ws = $wb.Worksheets.Item("Metadata")

# Update the Date property value
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new row for "Jurisdiction" (with an empty value) right before
# the "Description" row, shifting Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

Write-Output "done"
